$wb = $excel.ActiveWorkbook

# ALC!row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2083.3333
$ws.Cells.Item(40, 9).Value = 1933.3334
$ws.Cells.Item(40, 11).Value = 1933.3334
$ws.Cells.Item(40, 13).Value = -1758.3334

# ALC!row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 2616.125
$ws.Cells.Item(96, 9).Value = 1977.75
$ws.Cells.Item(96, 11).Value = 5933.25
$ws.Cells.Item(96, 13).Value = -4560.25

# ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16134.6
$ws.Cells.Item(32, 9).Value = 13549.571
$ws.Cells.Item(32, 11).Value = 13549.571
$ws.Cells.Item(32, 13).Value = -13262.571

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2550.4443
$ws.Cells.Item(132, 9).Value = 2550.4443
$ws.Cells.Item(132, 11).Value = 7651.3329
$ws.Cells.Item(132, 13).Value = -5121.3329

# BSM!row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1148.7142
$ws.Cells.Item(20, 9).Value = 1140.1666
$ws.Cells.Item(20, 11).Value = 1140.1666
$ws.Cells.Item(20, 13).Value = -893.1666

# BSM!row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 3438.8333
$ws.Cells.Item(94, 9).Value = 3626.6
$ws.Cells.Item(94, 11).Value = 3626.6
$ws.Cells.Item(94, 13).Value = -3175.6

# BSM!row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3131.75
$ws.Cells.Item(99, 9).Value = 2722
$ws.Cells.Item(99, 10).Value = 6000
$ws.Cells.Item(99, 11).Value = 2722
$ws.Cells.Item(99, 12).Value = 6000
$ws.Cells.Item(99, 13).Value = -1224
$ws.Cells.Item(99, 14).Value = -8996

# BSM!row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 867.2222
$ws.Cells.Item(107, 9).Value = 925.625
$ws.Cells.Item(107, 10).Value = 400
$ws.Cells.Item(107, 11).Value = 925.625
$ws.Cells.Item(107, 12).Value = 400
$ws.Cells.Item(107, 13).Value = 994.375
$ws.Cells.Item(107, 14).Value = -4240

# CRP!row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 13).ClearContents()

# CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 5862.857
$ws.Cells.Item(58, 9).Value = 3634.5
$ws.Cells.Item(58, 10).Value = 6754.2
$ws.Cells.Item(58, 11).Value = 3634.5
$ws.Cells.Item(58, 12).Value = 6754.2
$ws.Cells.Item(58, 13).Value = -3431.5
$ws.Cells.Item(58, 14).Value = -7160.2

# CRP!row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 33666.332
$ws.Cells.Item(68, 10).Value = 33666.332
$ws.Cells.Item(68, 12).Value = 33666.332
$ws.Cells.Item(68, 14).Value = -35164.332

# CRP!row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(71, 8).Value = 33666.332
$ws.Cells.Item(71, 10).Value = 33666.332
$ws.Cells.Item(71, 12).Value = 100998.996
$ws.Cells.Item(71, 14).Value = -108486.996

# CRP!row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 50250
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).ClearContents()

# CRP!row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(77, 8).Value = 50250
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).ClearContents()

# CRP!row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 10632.333
$ws.Cells.Item(86, 10).Value = 11666.333
$ws.Cells.Item(86, 12).Value = 11666.333
$ws.Cells.Item(86, 14).Value = -13912.333

# CRP!row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 10632.333
$ws.Cells.Item(89, 10).Value = 11666.333
$ws.Cells.Item(89, 12).Value = 58331.665
$ws.Cells.Item(89, 14).Value = -69563.66500000001

# CRP!row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 12648.667
$ws.Cells.Item(99, 9).Value = 8525.5
$ws.Cells.Item(99, 11).Value = 8525.5
$ws.Cells.Item(99, 13).Value = -7027.5

# CRP!row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 12648.667
$ws.Cells.Item(126, 9).Value = 8525.5
$ws.Cells.Item(126, 11).Value = 25576.5
$ws.Cells.Item(126, 13).Value = -23106.5

# CRP!row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3637.0667
$ws.Cells.Item(134, 9).Value = 2911.75
$ws.Cells.Item(134, 10).Value = 4466
$ws.Cells.Item(134, 11).Value = 8735.25
$ws.Cells.Item(134, 12).Value = 13398
$ws.Cells.Item(134, 13).Value = -6200.25
$ws.Cells.Item(134, 14).Value = -18468

# CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 5862.857
$ws.Cells.Item(136, 9).Value = 3634.5
$ws.Cells.Item(136, 10).Value = 6754.2
$ws.Cells.Item(136, 11).Value = 10903.5
$ws.Cells.Item(136, 12).Value = 20262.6
$ws.Cells.Item(136, 13).Value = -8353.5
$ws.Cells.Item(136, 14).Value = -25362.6

# CUL!row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 1087.5
$ws.Cells.Item(39, 9).Value = 1087.5
$ws.Cells.Item(39, 11).Value = 3262.5
$ws.Cells.Item(39, 13).Value = -2968.5

# CUL!row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()

# CUL!row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()

# GSM!row 26
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 70042
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 70042
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 70042
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(26, 14).Value = -70602

# GSM!row 50
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(50, 8).Value = 70042
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 70042
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 70042
$ws.Cells.Item(50, 13).ClearContents()
$ws.Cells.Item(50, 14).Value = -71038

# GSM!row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4977.25
$ws.Cells.Item(70, 9).Value = 4954.5
$ws.Cells.Item(70, 11).Value = 4954.5
$ws.Cells.Item(70, 13).Value = -4684.5

# GSM!row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 4977.25
$ws.Cells.Item(73, 9).Value = 4954.5
$ws.Cells.Item(73, 11).Value = 4954.5
$ws.Cells.Item(73, 13).Value = -4018.5

# LTW!row 44
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 14).ClearContents()

# LTW!row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3630.3684
$ws.Cells.Item(46, 9).Value = 2372.5
$ws.Cells.Item(46, 10).Value = 3965.8
$ws.Cells.Item(46, 11).Value = 2372.5
$ws.Cells.Item(46, 12).Value = 3965.8
$ws.Cells.Item(46, 13).Value = -2184.5
$ws.Cells.Item(46, 14).Value = -4341.8

# WVR!row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 15384.333
$ws.Cells.Item(45, 10).Value = 12508.571
$ws.Cells.Item(45, 12).Value = 12508.571
$ws.Cells.Item(45, 14).Value = -13490.571

# WVR!row 61
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(61, 8).Value = 70057
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 70057
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 70057
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(61, 14).Value = -70641

# WVR!row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 7704.7896
$ws.Cells.Item(62, 10).Value = 8055.9375
$ws.Cells.Item(62, 12).Value = 8055.9375
$ws.Cells.Item(62, 14).Value = -9303.9375

# WVR!row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 7704.7896
$ws.Cells.Item(65, 10).Value = 8055.9375
$ws.Cells.Item(65, 12).Value = 40279.6875
$ws.Cells.Item(65, 14).Value = -46519.6875
